$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: update label text (A3) and turn B3 into a hyperlink with Hyperlink style
$ws.Range("A3").Value = "Create any operation (receipt, distribution,transfer,adjustment)"
$ws.Range("B3").Value = "http://localhost:8080/openmrs/ws/rest/v2/inventory/stockOperation"
$ws.Hyperlinks.Add($ws.Range("B3"), "http://localhost:8080/openmrs/ws/rest/v2/inventory/stockOperation") | Out-Null
$ws.Range("B3").Style = "Hyperlink"

# Row 4 stays the same content (Complete Operation / stockOperation/{operation-id}) but keep as-is
$ws.Range("A4").Value = "Complete Operation"
$ws.Range("B4").Value = "http://localhost:8080/openmrs/ws/rest/v2/inventory/stockOperation/{operation-id}"

# Row 5 & 6 labels entered first (A5, A6), matching shared-string insertion order
$ws.Range("A5").Value = "Lab Item"
$ws.Range("A6").Value = "Pharmacy Item"

# Then the URLs: B6 (8093/...5452ec3e...) entered before B5 (8080/...2741bae2...)
$ws.Range("B6").Value = "http://localhost:8093/openmrs/ws/rest/v2/inventory/inventoryStockTakeSummary?limit=NaN&startIndex=1&stockroom_uuid=5452ec3e-2fe1-46de-8a6e-28c6442e4cc0"
$ws.Range("B5").Value = "http://localhost:8080/openmrs/ws/rest/v2/inventory/inventoryStockTakeSummary?limit=NaN&startIndex=1&stockroom_uuid=2741bae2-c5de-43ef-891f-7ec2fd58f442"

$ws.Hyperlinks.Add($ws.Range("B5"), "http://localhost:8080/openmrs/ws/rest/v2/inventory/inventoryStockTakeSummary?limit=NaN&startIndex=1&stockroom_uuid=2741bae2-c5de-43ef-891f-7ec2fd58f442") | Out-Null
$ws.Range("B5").Style = "Hyperlink"
$ws.Hyperlinks.Add($ws.Range("B6"), "http://localhost:8093/openmrs/ws/rest/v2/inventory/inventoryStockTakeSummary?limit=NaN&startIndex=1&stockroom_uuid=5452ec3e-2fe1-46de-8a6e-28c6442e4cc0") | Out-Null
$ws.Range("B6").Style = "Hyperlink"

# Adjust column widths to match the new layout
$ws.Columns.Item(1).ColumnWidth = 60.8776042
$ws.Columns.Item(2).ColumnWidth = 63.8776042

# Update the selection to match authored state
$ws.Range("B8").Select()
